# Insert two new data rows (Cereza / Lapins) right after the current row 310,
# pushing the existing rows 311-326 down to 313-328, then fill in the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 311 (each Insert() pushes everything below it down by one).
$ws.Rows.Item(311).Insert()
$ws.Rows.Item(311).Insert()

# --- New row 311: Lapins / Primera ---
$ws.Cells.Item(311, 1).Value  = 9
$ws.Cells.Item(311, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(311, 3).Value  = "Metropolitana"
$ws.Cells.Item(311, 4).Value  = 44610
$ws.Cells.Item(311, 5).Value  = 13
$ws.Cells.Item(311, 6).Value  = "Fruta"
$ws.Cells.Item(311, 7).Value  = 100103
$ws.Cells.Item(311, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(311, 9).Value  = 100103001
$ws.Cells.Item(311, 10).Value = "Cereza"
$ws.Cells.Item(311, 11).Value = "Lapins"
$ws.Cells.Item(311, 12).Value = "Primera"
$ws.Cells.Item(311, 13).Value = 180
$ws.Cells.Item(311, 14).Value = 5000
$ws.Cells.Item(311, 15).Value = 5000
$ws.Cells.Item(311, 16).Value = 5000
$ws.Cells.Item(311, 17).Value = "$/bandeja 5 kilos"
$ws.Cells.Item(311, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(311, 19).Value = 1000
$ws.Cells.Item(311, 20).Value = 5

# --- New row 312: Lapins / Segunda ---
$ws.Cells.Item(312, 1).Value  = 9
$ws.Cells.Item(312, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(312, 3).Value  = "Metropolitana"
$ws.Cells.Item(312, 4).Value  = 44610
$ws.Cells.Item(312, 5).Value  = 13
$ws.Cells.Item(312, 6).Value  = "Fruta"
$ws.Cells.Item(312, 7).Value  = 100103
$ws.Cells.Item(312, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(312, 9).Value  = 100103001
$ws.Cells.Item(312, 10).Value = "Cereza"
$ws.Cells.Item(312, 11).Value = "Lapins"
$ws.Cells.Item(312, 12).Value = "Segunda"
$ws.Cells.Item(312, 13).Value = 260
$ws.Cells.Item(312, 14).Value = 4000
$ws.Cells.Item(312, 15).Value = 4000
$ws.Cells.Item(312, 16).Value = 4000
$ws.Cells.Item(312, 17).Value = "$/bandeja 5 kilos"
$ws.Cells.Item(312, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(312, 19).Value = 800
$ws.Cells.Item(312, 20).Value = 5
